# Apply the cryptos list update (prices + 1h volume deltas) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.129.03'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '2.428.86'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''563.02'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = '''144.14'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '2.428.63'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("E10").Value = '  +0.18%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("E12").Value = '  -2.74%  '
$ws.Range("D13").Value = '''0.350'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").Value = '''26.30'
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '62.009.66'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = '2.431.27'
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").Value = '''11.26'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").Value = '''323.39'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '''67.57'
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  -3.87%  '
$ws.Range("D27").Value = '''554.44'
$ws.Range("E27").Value = '  -4.52%  '
$ws.Range("D28").Value = '2.546.13'
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").Value = '''0.996'
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").Value = '0.0₃0939'
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("D31").Value = '''8.26'
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").Value = '''1.52'
$ws.Range("E35").Value = '  -1.78%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = '''4.78'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("E39").Value = '  -3.48%  '
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("E42").Value = '  -1.15%  '
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("D44").Value = '''2.28'
$ws.Range("E44").Value = '  -2.50%  '
$ws.Range("D45").Value = '''147.79'
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("E47").Value = '  -1.39%  '
$ws.Range("D48").Value = '''20.05'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("E50").Value = '  +0.17%  '
$ws.Range("E51").Value = '  +0.22%  '
